# Add a new "SI Table" worksheet after Sheet1
$wb = $excel.ActiveWorkbook

$siTable = $wb.Worksheets.Add()
$siTable.Name = "SI Table"

$ws1 = $wb.Worksheets.Item("Sheet1")
$siTable.Move($null, $ws1)

# Moving the new sheet shifts worksheet positions, so re-resolve the
# "Sheet1" reference before editing it.
$ws1 = $wb.Worksheets.Item("Sheet1")

# Append the new "Table 2" block to Sheet1 starting at row 23
$ws1.Range("A23").Value = "Table 2. Mean concentration values of potential hydrologic ﬂowpath tracers in each end-member. Units for Cl_x0001_.."

$ws1.Range("A24").Value = "Hungerford"

$ws1.Range("A25").Value = "End-member"
$ws1.Range("B25").Value = "Cl-"
$ws1.Range("C25").Value = "Cl- STDV"
$ws1.Range("D25").Value = "Na+"
$ws1.Range("E25").Value = "Ca2+"
$ws1.Range("F25").Value = "Mg2+"
$ws1.Range("G25").Value = "H4SiO4"
$ws1.Range("H25").Value = "d18O"
$ws1.Range("I25").Value = "dD"
$ws1.Range("J25").Value = "TOC"

$ws1.Range("A26").Value = "Snowmelt"
$ws1.Range("B26").Formula = "=AVERAGE(U18:U20)"
$ws1.Range("C26").Formula = "=STDEV(U18:U20)"

$ws1.Range("A27").Value = "Soil water wet site"
$ws1.Range("A28").Value = "Soil water dry site"
$ws1.Range("A29").Value = "Groundwater/baseflow"

$ws1.Range("A30").Value = "Wade"
$ws1.Range("A31").Value = "Snowmelt"
$ws1.Range("A32").Value = "Soil water wet site"
$ws1.Range("A33").Value = "Soil water dry site"
$ws1.Range("A34").Value = "Groundwater/baseflow"

$ws1.Range("I24").Select() | Out-Null
